# Changes regarding the MTF tiers
$wb = $excel.ActiveWorkbook

# Worksheets (by name, to be safe regardless of index assumptions)
$wsInfo = $wb.Worksheets.Item("ScenarioInfo")
$wsParams = $wb.Worksheets.Item("ScenarioParameters")
$wsData = $wb.Worksheets.Item("SpecsData")

# Update ScenarioParameters values
$wsParams.Range("C2").Value = 4

# Update SpecsData values (MTF tier related figures)
$wsData.Range("M2").Value = 4.3
$wsData.Range("N2").Value = 4.5
$wsData.Range("O2").Value = 0.83
$wsData.Range("Q2").Value = 0.076

# Update selections on each sheet to reflect new active cell per sheet
$wsInfo.Range("C1").Select() | Out-Null
$wsParams.Range("E7").Select() | Out-Null
$wsData.Range("R2").Select() | Out-Null

# Make ScenarioParameters the active (selected) sheet/tab
$wsParams.Activate()
$wsParams.Range("E7").Select() | Out-Null
